$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain text, preserving exact formatting
# (values like "317.30", "1.00", "0.110" would otherwise be auto-converted to numbers)
$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.465.54'
$ws.Range("D2").Style = $origStyle
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.336.24'
$ws.Range("D3").Style = $origStyle
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.30'
$ws.Range("D5").Style = $origStyle
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.27'
$ws.Range("D6").Style = $origStyle
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.634'
$ws.Range("D7").Style = $origStyle
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.39'
$ws.Range("D10").Style = $origStyle
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.24'
$ws.Range("D12").Style = $origStyle
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.971'
$ws.Range("D14").Style = $origStyle
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.77'
$ws.Range("D15").Style = $origStyle
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.688.92'
$ws.Range("D16").Style = $origStyle
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.339.05'
$ws.Range("D17").Style = $origStyle
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.449.61'
$ws.Range("D18").Style = $origStyle
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.61'
$ws.Range("D19").Style = $origStyle
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000104'
$ws.Range("D20").Style = $origStyle
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.72'
$ws.Range("D21").Style = $origStyle
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '259.40'
$ws.Range("D23").Style = $origStyle
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.48'
$ws.Range("D25").Style = $origStyle
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = $origStyle
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.21'
$ws.Range("D27").Style = $origStyle
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.95'
$ws.Range("D28").Style = $origStyle
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.06'
$ws.Range("D30").Style = $origStyle
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.25'
$ws.Range("D31").Style = $origStyle
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0881'
$ws.Range("D32").Style = $origStyle
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.94'
$ws.Range("D33").Style = $origStyle
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.93'
$ws.Range("D34").Style = $origStyle
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.110'
$ws.Range("D36").Style = $origStyle
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.51'
$ws.Range("D37").Style = $origStyle
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.70'
$ws.Range("D39").Style = $origStyle
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.62'
$ws.Range("D40").Style = $origStyle
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.45'
$ws.Range("D41").Style = $origStyle
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '68.25'
$ws.Range("D43").Style = $origStyle
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '114.31'
$ws.Range("D45").Style = $origStyle
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.46'
$ws.Range("D46").Style = $origStyle
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.08'
$ws.Range("D47").Style = $origStyle
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.47'
$ws.Range("D48").Style = $origStyle
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.01'
$ws.Range("D49").Style = $origStyle
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.541.68'
$ws.Range("D50").Style = $origStyle

# Other cells (Coin names, Links, Volume percentages) - plain text, safe to set directly
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("B32").Value = 'Hedera'
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("B48").Value = 'Celestia'
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("C48").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E2").Value = '  -2.55%  '
$ws.Range("E3").Value = '  -3.30%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("E5").Value = '  -2.50%  '
$ws.Range("E6").Value = '  -2.18%  '
$ws.Range("E7").Value = '  -1.81%  '
$ws.Range("E9").Value = '  -6.68%  '
$ws.Range("E10").Value = '  -7.32%  '
$ws.Range("E11").Value = '  -3.52%  '
$ws.Range("E12").Value = '  -6.58%  '
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("E14").Value = '  -6.24%  '
$ws.Range("E15").Value = '  -10.14%  '
$ws.Range("E16").Value = '  -3.50%  '
$ws.Range("E17").Value = '  -3.00%  '
$ws.Range("E18").Value = '  -2.67%  '
$ws.Range("E19").Value = '  +2.01%  '
$ws.Range("E20").Value = '  -5.01%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("E24").Value = '  -7.82%  '
$ws.Range("E25").Value = '  -4.99%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  -7.34%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("E29").Value = '  +1.54%  '
$ws.Range("E30").Value = '  -2.85%  '
$ws.Range("E31").Value = '  -8.20%  '
$ws.Range("E32").Value = '  -6.17%  '
$ws.Range("E33").Value = '  -9.26%  '
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("E35").Value = '  -2.49%  '
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("E37").Value = '  -8.49%  '
$ws.Range("E38").Value = '  -6.08%  '
$ws.Range("E39").Value = '  -7.73%  '
$ws.Range("E40").Value = '  -9.91%  '
$ws.Range("E41").Value = '  -10.94%  '
$ws.Range("E42").Value = '  -3.00%  '
$ws.Range("E43").Value = '  -2.89%  '
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("E45").Value = '  -7.14%  '
$ws.Range("E46").Value = '  -4.67%  '
$ws.Range("E47").Value = '  +10.23%  '
$ws.Range("E48").Value = '  -9.92%  '
$ws.Range("E49").Value = '  -4.74%  '
$ws.Range("E50").Value = '  -2.40%  '
$ws.Range("E51").Value = '  -3.00%  '
